$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("1-torta,", "10-10-2019", "10-10-2019", "lal",       "ii@ff.com",        "999999999", "Pendiente", [double]30000, [double]0, [double]0),
    @("1-torta,", "10-10-2019", "13-10-2019", "lala lala", "lala@gmail.com",   "999999999", "Pendiente", [double]30000, [double]0, [double]0),
    @("1-torta,", "10-10-2019", "15-10-2019", "lala",      "isa@gmail.com",    "999999999", "Pendiente", [double]30000, [double]0, [double]0)
)

$startRow = 7
# Columns (1-based) whose literal text would otherwise be auto-detected by
# Excel as a date/number (the order-date, pickup-date and phone columns);
# those need Text format applied first so the literal string is preserved.
$textColumns = @(2, 3, 6)

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]
    for ($c = 0; $c -lt $row.Count; $c++) {
        $cell = $ws.Cells.Item($r, $c + 1)
        if ($textColumns -contains ($c + 1)) {
            $cell.NumberFormat = "@"
        }
        $cell.Value = $row[$c]
    }
}
